$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Fale"
$ws.Range("A3").Value = "Ruso"
$ws.Range("A4").Value = "Lope"
$ws.Range("A5").Value = "Puche"
$ws.Range("A7").Value = "Coquina"
